# feat: add 2022-Q1 data
#
# 1. Create a new "2022-Q1" sheet (positioned right after "2021-Q4" and
#    before "总计") by duplicating the "2021-Q4" sheet so it inherits the
#    same layout/styles, then overwrite its contents with the new quarter's
#    fund-holding data (2 funds instead of 4).
# 2. Prepend a new "2022-Q1" row to the "总计" (totals) summary sheet and
#    renumber the existing rows' index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: new "2022-Q1" worksheet
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $q4)
$newWs = $wb.Worksheets.Item($q4.Index + 1)
$newWs.Name = "2022-Q1"

# the template ("2021-Q4") has 4 data rows (rows 2-5); the new quarter only
# has 2 funds, so drop the extra two template rows.
$newWs.Range("A4:A5").EntireRow.Delete()

# Row 2: 006601 / 国融融泰灵活配置混合A
$newWs.Range("B2").Value = "'006601"
$newWs.Range("B2").Style = "Normal"
$newWs.Range("C2").Value = "国融融泰灵活配置混合A"
$newWs.Range("C2").Style = "Normal"
$newWs.Range("D2").Value = "'0.04"
$newWs.Range("D2").Style = "Normal"
$newWs.Range("E2").Value = "'47.44"
$newWs.Range("E2").Style = "Normal"
$newWs.Range("F2").Value = "'3.90"
$newWs.Range("F2").Style = "Normal"
$newWs.Range("G2").Value = "'0.0016"
$newWs.Range("G2").Style = "Normal"
$newWs.Range("H2").Value = 3

# Row 3: 006602 / 国融融泰灵活配置混合C
$newWs.Range("B3").Value = "'006602"
$newWs.Range("B3").Style = "Normal"
$newWs.Range("C3").Value = "国融融泰灵活配置混合C"
$newWs.Range("C3").Style = "Normal"
$newWs.Range("D3").Value = "'0.01"
$newWs.Range("D3").Style = "Normal"
$newWs.Range("E3").Value = "'47.44"
$newWs.Range("E3").Style = "Normal"
$newWs.Range("F3").Value = "'3.90"
$newWs.Range("F3").Style = "Normal"
$newWs.Range("G3").Value = "'0.0004"
$newWs.Range("G3").Style = "Normal"
$newWs.Range("H3").Value = 3

# ---------------------------------------------------------------------
# Step 2: update the "总计" summary sheet
# ---------------------------------------------------------------------
$totalWs = $wb.Worksheets.Item("总计")

# insert a new row under the header and seed it with the row directly
# below (copies the existing number formatting/border style), then
# overwrite with the 2022-Q1 totals.
$totalWs.Rows.Item(2).Insert()
$totalWs.Range("A3:D3").Copy($totalWs.Range("A2:D2"))

$totalWs.Range("B2").Value = "2022-Q1"
$totalWs.Range("C2").Value = 2
$totalWs.Range("D2").Value = 0

# renumber the index column (A) for the rows that shifted down
$totalWs.Range("A3").Value = 1
$totalWs.Range("A4").Value = 2
$totalWs.Range("A5").Value = 3
$totalWs.Range("A6").Value = 4
$totalWs.Range("A7").Value = 5

# restore the originally-active sheet/tab selection (copying/renaming
# sheets moves the active tab, which this edit should not change).
$wb.Worksheets.Item(1).Activate()
